$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grade average (not ortalaması)
$ws.Range("H3").Formula = "=AVERAGE(D4:D12)"

# Student number (Numara:)
$ws.Range("L4").Value = 20215070055

# Highest age (en büyük yaş)
$ws.Range("H5").Formula = "=MAX(E4:E12)"

# Name Surname (Ad Soyad:)
$ws.Range("L5").Value = "Muhammed Ali Harmancı"

# Department (Bölüm:)
$ws.Range("L6").Value = "Yönetim Bilişim Sistemleri"

# Lowest grade (en küçük not)
$ws.Range("H8").Formula = "=MIN(D4:D12)"

# 2nd lowest age (en küçük 2. yaş)
$ws.Range("H10").Formula = "=SMALL(E4:E12,2)"

# 4th highest grade (en büyük 4. not)
$ws.Range("H11").Formula = "=LARGE(D4:D12,4)"

# Today's date (bugünün tarihi)
$ws.Range("H13").Formula = "=DATE(2021,12,13)"

# Today's date and time (bugün tarih ve saati)
$ws.Range("H14").Formula = "=NOW()"

# Move the active selection to K6 like in the edited workbook
$ws.Range("K6").Select() | Out-Null
